$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 198, pushing existing rows 198:232 down to 199:233
$ws.Rows("198:198").Insert()

# Populate the newly inserted row 198 with the new weekly price record
$ws.Range("A198").Value = 8
$ws.Range("B198").Value = "Terminal La Palmera de La Serena"
$ws.Range("C198").Value = "Coquimbo"
$ws.Range("D198").Value = 44522
$ws.Range("E198").Value = 4
$ws.Range("F198").Value = 100114013
$ws.Range("G198").Value = "Zanahoria"
$ws.Range("H198").Value = "Sin especificar"
$ws.Range("I198").Value = "Primera"
$ws.Range("J198").Value = 600
$ws.Range("K198").Value = 6500
$ws.Range("L198").Value = 7000
$ws.Range("M198").Value = 6750
$ws.Range("N198").Value = "$/saco 20 kilos"
$ws.Range("O198").Value = "Provincia del Elquí"
$ws.Range("P198").Value = 338
$ws.Range("Q198").Value = 20
$ws.Range("R198").Value = "Hortaliza"
